$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 18 with url / field_link_api / FALSE
$ws.Range("A18").Value = "url"
$ws.Range("B18").Value = "field_link_api"
$ws.Range("C18").Value = $false

# Update the selection to reflect where the cursor moved after the edit
$ws.Range("A19").Select()
